# Rename the three header/footer logo pictures (wp:docPr / pic:cNvPr
# "name" attribute) so the Pearson logos swap from image1.png to
# image2.png and the BTEC logo swaps from image2.jpg to image1.jpg.
#
# InlineShape has no .Name property in Word's object model (that only
# exists on Shape), so each picture is briefly converted to a floating
# Shape, renamed, then converted back to an inline picture.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($rng, $index, $newName) {
    $inlineShape = $rng.InlineShapes.Item($index)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Footer, primary/default (type="default" -> footer2.xml, PearsonLogo id=2)
Rename-InlinePicture $sec.Footers.Item(1).Range 1 "image2.png"

# Footer, first page (type="first" -> footer1.xml, PearsonLogo id=3)
Rename-InlinePicture $sec.Footers.Item(2).Range 1 "image2.png"

# Header, first page (type="first" -> header1.xml, BTec_Logo-Orange id=1)
Rename-InlinePicture $sec.Headers.Item(2).Range 1 "image1.jpg"
